# Generate Report for Handback
# -----------------------------------------------------------------------
# This mirrors the localization-report generator's "handback" pass:
#   - flips each row's Status from "Ready for handoff" to
#     "Handed back: in sync with en-US" (shows up everywhere that text
#     is used: the Overview rollup and each language sheet)
#   - records the freshly produced "Latest Target File" / "Latest
#     Handback File" links (mirrors of the source + handoff-target
#     filenames) for every row
#   - stamps "Latest Handback DateTime" with the actual handback time
#     per language
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/6ff6d209de10b2743163637e9ca7a45e0973bbc3/e2e/a.md"
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01687cf634aadee498224d38c146477e921f43f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f2d935c9c0c65a1bc1cb723e1e24b3d3c65a2bf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$zhHandbackTime = "2016-03-24 11:23:50"
$deHandbackTime = "2016-03-24 11:24:00"

# ---- Overview rollup: Status shown for both languages on each row ----
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $aMdUrl, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, "", "", $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $aMdUrl, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, "", "", $zhXlfName)

$wsZh.Range("H2").Value = $zhHandbackTime
$wsZh.Range("H3").Value = $zhHandbackTime

# ---- de-de sheet ----
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $aMdUrl, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, "", "", $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $aMdUrl, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, "", "", $deXlfName)

$wsDe.Range("H2").Value = $deHandbackTime
$wsDe.Range("H3").Value = $deHandbackTime
